# 🔄 Actualización automática del tracker
# Append new match rows (164-170) to the results tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 164
$ws.Cells.Item(164, 1).Value = 14851756
$ws.Cells.Item(164, 2).NumberFormat = "@"
$ws.Cells.Item(164, 2).Value = "2025-10-15"
$ws.Cells.Item(164, 2).ClearFormats()
$ws.Cells.Item(164, 3).Value = "Marcos Giron"
$ws.Cells.Item(164, 4).Value = "Alejandro Davidovich Fokina"
$ws.Cells.Item(164, 5).Value = "Gana Alejandro Davidovich Fokina"
$ws.Cells.Item(164, 6).Value = 1.62

# Row 165
$ws.Cells.Item(165, 1).Value = 14851812
$ws.Cells.Item(165, 2).NumberFormat = "@"
$ws.Cells.Item(165, 2).Value = "2025-10-15"
$ws.Cells.Item(165, 2).ClearFormats()
$ws.Cells.Item(165, 3).Value = "Marin Čilić"
$ws.Cells.Item(165, 4).Value = "Casper Ruud"
$ws.Cells.Item(165, 5).Value = "Gana Casper Ruud"
$ws.Cells.Item(165, 6).Value = 1.44

# Row 166
$ws.Cells.Item(166, 1).Value = 14859066
$ws.Cells.Item(166, 2).NumberFormat = "@"
$ws.Cells.Item(166, 2).Value = "2025-10-16"
$ws.Cells.Item(166, 2).ClearFormats()
$ws.Cells.Item(166, 3).Value = "McCartney Kessler"
$ws.Cells.Item(166, 4).Value = "Liudmila Samsonova"
$ws.Cells.Item(166, 5).Value = "Gana Liudmila Samsonova"
$ws.Cells.Item(166, 6).Value = 1.91

# Row 167
$ws.Cells.Item(167, 1).Value = 14859094
$ws.Cells.Item(167, 2).NumberFormat = "@"
$ws.Cells.Item(167, 2).Value = "2025-10-16"
$ws.Cells.Item(167, 2).ClearFormats()
$ws.Cells.Item(167, 3).Value = "Elise Mertens"
$ws.Cells.Item(167, 4).Value = "Tereza Valentova"
$ws.Cells.Item(167, 5).Value = "Gana Elise Mertens"
$ws.Cells.Item(167, 6).Value = 2

# Row 168
$ws.Cells.Item(168, 1).Value = 14858866
$ws.Cells.Item(168, 2).NumberFormat = "@"
$ws.Cells.Item(168, 2).Value = "2025-10-16"
$ws.Cells.Item(168, 2).ClearFormats()
$ws.Cells.Item(168, 3).Value = "Inaki Montes-de la Torre"
$ws.Cells.Item(168, 4).Value = "Ioannis Xilas"
$ws.Cells.Item(168, 5).Value = "Gana Ioannis Xilas"
$ws.Cells.Item(168, 6).Value = 2.63

# Row 169
$ws.Cells.Item(169, 1).Value = 14857967
$ws.Cells.Item(169, 2).NumberFormat = "@"
$ws.Cells.Item(169, 2).Value = "2025-10-15"
$ws.Cells.Item(169, 2).ClearFormats()
$ws.Cells.Item(169, 3).Value = "Dhakshineswar Suresh"
$ws.Cells.Item(169, 4).Value = "Alfredo Perez"
$ws.Cells.Item(169, 5).Value = "Gana Dhakshineswar Suresh"
$ws.Cells.Item(169, 6).Value = 1.73

# Row 170
$ws.Cells.Item(170, 1).Value = 14863015
$ws.Cells.Item(170, 2).NumberFormat = "@"
$ws.Cells.Item(170, 2).Value = "2025-10-15"
$ws.Cells.Item(170, 2).ClearFormats()
$ws.Cells.Item(170, 3).Value = "Harriet Dart"
$ws.Cells.Item(170, 4).Value = "Dasha Plekhanova"
$ws.Cells.Item(170, 5).Value = "Gana Dasha Plekhanova"
$ws.Cells.Item(170, 6).Value = 3.75
